$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Header rich-text strings: "Volume 31 Number 10" -> "...11"
# and "Report Covering the Week 3/4/2024 Through 3/10/2024"
# -> "...3/11/2024 Through 3/17/2024"
# ---------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "11"
$ws.Range("C9").Characters(27, 8).Text = "3/11/2024"
$ws.Range("C9").Characters(47, 9).Text = "3/17/2024"

# ---------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------
$ws.Range("D15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -75
$ws.Range("J15").Value = 6
$ws.Range("K15").Value = -66.666666666666
$ws.Range("N15").Value = -90.909090909090

# ---------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -25
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 38
$ws.Range("J16").Value = 36
$ws.Range("K16").Value = 5.555555555555
$ws.Range("L16").Value = -33.333333333333
$ws.Range("M16").Value = -48.648648648648
$ws.Range("N16").Value = -91.574279379157

# ---------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -75
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = -61.764705882352
$ws.Range("I17").Value = 55
$ws.Range("J17").Value = 80
$ws.Range("K17").Value = -31.25
$ws.Range("L17").Value = -23.611111111111
$ws.Range("M17").Value = -28.571428571428
$ws.Range("N17").Value = -66.257668711656

# ---------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -28.571428571428
$ws.Range("I18").Value = 23
$ws.Range("J18").Value = 36
$ws.Range("K18").Value = -36.111111111111
$ws.Range("L18").Value = -42.5
$ws.Range("M18").Value = -53.061224489795
$ws.Range("N18").Value = -96.587537091988

# ---------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 142.857142857143
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = 21.621621621621
$ws.Range("I19").Value = 98
$ws.Range("J19").Value = 125
$ws.Range("K19").Value = -21.6
$ws.Range("L19").Value = -16.949152542372
$ws.Range("M19").Value = -7.547169811320
$ws.Range("N19").Value = -51.485148514851

# ---------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 7
$ws.Range("H20").Value = 40
$ws.Range("I20").Value = 25
$ws.Range("J20").Value = 20
$ws.Range("K20").Value = 25
$ws.Range("L20").Value = -34.210526315789
$ws.Range("M20").Value = -44.444444444444
$ws.Range("N20").Value = -94.791666666666

# ---------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = 8.333333333333
$ws.Range("F21").Value = 88
$ws.Range("G21").Value = 110
$ws.Range("H21").Value = -20
$ws.Range("I21").Value = 241
$ws.Range("J21").Value = 304
$ws.Range("K21").Value = -20.723684210526
$ws.Range("L21").Value = -28.273809523809
$ws.Range("M21").Value = -32.303370786516
$ws.Range("N21").Value = -87.956021989005

# ---------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 5
$ws.Range("K22").Value = 25
$ws.Range("L22").Value = 66.666666666666
$ws.Range("M22").Value = -16.666666666666

# ---------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 54
$ws.Range("E24").Value = -62.962962962963
$ws.Range("F24").Value = 160
$ws.Range("G24").Value = 153
$ws.Range("H24").Value = 4.575163398692
$ws.Range("I24").Value = 421
$ws.Range("J24").Value = 360
$ws.Range("K24").Value = 16.944444444444
$ws.Range("L24").Value = 50.357142857142
$ws.Range("M24").Value = 76.150627615062

# ---------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 40
$ws.Range("E25").Value = -72.5
$ws.Range("F25").Value = 84
$ws.Range("G25").Value = 101
$ws.Range("H25").Value = -16.831683168316
$ws.Range("I25").Value = 245
$ws.Range("J25").Value = 211
$ws.Range("K25").Value = 16.113744075829
$ws.Range("L25").Value = 78.832116788321

# ---------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------
$ws.Range("C26").Value = 17
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = 30.769230769230
$ws.Range("F26").Value = 58
$ws.Range("G26").Value = 42
$ws.Range("H26").Value = 38.095238095238
$ws.Range("I26").Value = 136
$ws.Range("J26").Value = 109
$ws.Range("K26").Value = 24.770642201834
$ws.Range("L26").Value = 24.770642201834
$ws.Range("M26").Value = -13.375796178343

# ---------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------
$ws.Range("D27").Value = 4
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 11
$ws.Range("K27").Value = -45.454545454545
$ws.Range("L27").Value = -62.5

# ---------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 2
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 233.333333333333
$ws.Range("I28").Value = 16
$ws.Range("J28").Value = 8
$ws.Range("K28").Value = 100
$ws.Range("L28").Value = 60

# ---------------------------------------------------------------
# Row 29 - Shooting Vic. (text "0"/"***.*" cells become real numbers)
# ---------------------------------------------------------------
$ws.Range("D29").Value = 1
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("E29").Value = -100
$ws.Range("E29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G29").Value = 1
$ws.Range("G29").NumberFormat = "#,##0"
$ws.Range("H29").Value = -100
$ws.Range("H29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J29").Value = 2

# ---------------------------------------------------------------
# Row 30 - Shooting Inc. (text "0"/"***.*" cells become real numbers)
# ---------------------------------------------------------------
$ws.Range("D30").Value = 1
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G30").Value = 1
$ws.Range("G30").NumberFormat = "#,##0"
$ws.Range("H30").Value = -100
$ws.Range("H30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J30").Value = 2

# ---------------------------------------------------------------
# Row 31 - Hate Crimes (numeric cells revert to text placeholders)
# ---------------------------------------------------------------
$ws.Range("D31").Value = "'0"
$ws.Range("E31").Value = "'***.*"
$ws.Range("C15").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("N23").Copy()
$ws.Range("E31").PasteSpecial(-4122)

# ---------------------------------------------------------------
# Row 33 - Traffic Fatalities (numeric cells revert to text placeholders)
# ---------------------------------------------------------------
$ws.Range("C33").Value = "'0"
$ws.Range("G33").Value = "'0"
$ws.Range("H33").Value = "'***.*"
$ws.Range("C15").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("G33").PasteSpecial(-4122)
$ws.Range("N23").Copy()
$ws.Range("H33").PasteSpecial(-4122)

Write-Host "Edit complete"
